# Card12: add a new service-log event in row 20 (below the existing last
# row, 19), matching the same "card number / date / serviced-by / event /
# correction" pattern used by the other manually-logged rows (14-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

$row = 20

# Column A holds the card number as text (same as every other row in this
# sheet, e.g. A19 = "12"). Excel would normally infer a plain "12" as a
# number, so force text formatting first, write the value, then restore
# the default "Normal" style so no stray number-format style is left
# behind on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12"
$ws.Cells.Item($row, 1).Style = "Normal"

# Columns B-K (Tones / element checkmarks) are left blank for this event,
# same as rows 14-19.

# L = Date, M = Serviced by, O = Correction (N = Event is left blank).
$ws.Cells.Item($row, 12).Value = "14\8\2024"
$ws.Cells.Item($row, 13).Value = "تيم العمل"
$ws.Cells.Item($row, 15).Value = "تم تشحيم المكنه بالكامل +عمل صيانه"
